$wb = $excel.ActiveWorkbook

# Default workbook font: Arial -> Calibri
$normalStyle = $wb.Styles.Item("Normal")
$normalStyle.Font.Name = "Calibri"

# Add "size" sheet after "san pham"
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetSize = $wb.Worksheets.Add([Type]::Missing, $afterSheet)
$sheetSize.Name = "size"
$sheetSize.Range("A1").Value = "ten"
$sheetSize.Range("B1").Value = "trangthai"
$sheetSize.Range("A2").Value = "XXL"
$sheetSize.Range("B2").Value = 1
[void]$sheetSize.Range("A3").Select()

# Add "mausac" sheet after "size"
$sheetMauSac = $wb.Worksheets.Add([Type]::Missing, $sheetSize)
$sheetMauSac.Name = "mausac"
$sheetMauSac.Range("A1").Value = "ten"
$sheetMauSac.Range("B1").Value = "trangthai"
$sheetMauSac.Range("A2").Value = "Kem"
$sheetMauSac.Range("B2").Value = 1
[void]$sheetMauSac.Range("A3").Select()
